# feats. SPI & heartbeat
# Slide 3 ("Bobby-MyRPC" architecture figure): replace the two empty
# rounded-rectangle placeholders + the two "BeanPostProcessor" textboxes
# with a small "Client / ClientProxy / 序列化 / 协议编码 / 网络传输 ..." /
# "Server / ... / 反序列化 / 协议解码 / 网络传输" call-chain diagram made of
# nine rounded rectangles.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$EMU_PER_PT = 12700
# Tiny nudge so the Single-precision round trip inside Shape.Left/Top/
# Width/Height lands on the exact target EMU instead of one EMU short.
$EPS = 0.00003

function EmuToPt($emu) {
    return ($emu / $EMU_PER_PT) + $EPS
}

# ---------------------------------------------------------------------
# 1. Drop the two "BeanPostProcessor" textboxes entirely.
# ---------------------------------------------------------------------
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "文本框 5" -or $sh.Name -eq "文本框 6") {
        $sh.Delete()
    }
}

# ---------------------------------------------------------------------
# 2. Identify the two surviving rounded rectangles ("矩形: 圆角 1" /
#    "矩形: 圆角 4"). We keep "矩形: 圆角 4" (becomes the "Client" box) and
#    use it as the style template for every new box; "矩形: 圆角 1" gets
#    removed once we no longer need it.
# ---------------------------------------------------------------------
$template = $null
$toRemove = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "矩形: 圆角 4") {
        $template = $sh
    } elseif ($sh.Name -eq "矩形: 圆角 1") {
        $toRemove = $sh
    }
}

# ---------------------------------------------------------------------
# 3. Burn through the runtime's low "free id" slots with throw-away
#    duplicates of the template so the *real* new shapes line up on the
#    ids the final deck actually uses (8,9,10,...).
# ---------------------------------------------------------------------
$junk1 = $template.Duplicate().Item(1)
$junk2 = $template.Duplicate().Item(1)
$junk1.Delete()
$junk2.Delete()

# The empty "矩形: 圆角 1" box is no longer needed.
$toRemove.Delete()

# ---------------------------------------------------------------------
# 4. Reposition/resize/retext the kept box -> "Client".
# ---------------------------------------------------------------------
$template.Left = EmuToPt 1436116
$template.Top = EmuToPt 1274923
$template.Width = EmuToPt 1654499
$template.Height = EmuToPt 725894
$template.TextFrame.TextRange.Text = "Client"

# ---------------------------------------------------------------------
# 5. Add the remaining eight boxes by duplicating the template (keeps
#    the rounded-rect fill/line/quick-style), then moving/retexting/
#    renaming each one.
# ---------------------------------------------------------------------
$newBoxes = @(
    @{ name = "矩形: 圆角 7";  x = 7715837; y = 1274923; cx = 1654499; cy = 725894; text = "Server" }
    @{ name = "矩形: 圆角 8";  x = 1436116; y = 2370392; cx = 1654499; cy = 725894; text = "ClientProxy" }
    @{ name = "矩形: 圆角 9";  x = 1436116; y = 4648147; cx = 1654499; cy = 418073; text = "序列化" }
    @{ name = "矩形: 圆角 10"; x = 1436116; y = 5245675; cx = 1654499; cy = 418073; text = "协议编码" }
    @{ name = "矩形: 圆角 11"; x = 1436116; y = 5843203; cx = 1654499; cy = 418073; text = "网络传输" }
    @{ name = "矩形: 圆角 12"; x = 7637739; y = 4648147; cx = 1654499; cy = 418073; text = "反序列化" }
    @{ name = "矩形: 圆角 13"; x = 7637739; y = 5245675; cx = 1654499; cy = 418073; text = "协议解码" }
    @{ name = "矩形: 圆角 14"; x = 7637739; y = 5843203; cx = 1654499; cy = 418073; text = "网络传输" }
)

foreach ($box in $newBoxes) {
    $ns = $template.Duplicate().Item(1)
    $ns.Name = $box.name
    $ns.Left = EmuToPt $box.x
    $ns.Top = EmuToPt $box.y
    $ns.Width = EmuToPt $box.cx
    $ns.Height = EmuToPt $box.cy
    $ns.TextFrame.TextRange.Text = $box.text
}
